# Add a new "Folio No*" column header in H1 (new shared string) and
# move the active selection to H2, matching the author's "dynamic
# dropdowns for commitment" edit that introduced a Folio No column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Folio No*"

$ws.Range("H2").Select()
